$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the text in B10: "Coagulació, floculació" -> "Coagulació, floculació i decantació"
$ws.Range("B10").Value = "Coagulació, floculació i decantació"

# Update A14: "MF_UF" -> "MF/UF"
$ws.Range("A14").Value = "MF/UF"

# Update B14: "Microfiltració, ultrafiltració" -> "Microfiltració o ultrafiltració"
$ws.Range("B14").Value = "Microfiltració o ultrafiltració"

# Update the active cell selection to A15
$ws.Range("A15").Select()
